$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated win/loss/draw probability values from games pulled March 7
$updates = @{
    "B2" = 0.2256532066508314
    "C2" = 0.498812351543943
    "J2" = 0.01425178147268409
    "P2" = 0.1425178147268409
    "S2" = 0.1187648456057007
    "B3" = 0.00909090909090909
    "C3" = 0.02727272727272727
    "J3" = 0.03181818181818181
    "P3" = 0.6863636363636364
    "S3" = 0.2454545454545455
    "J4" = 0.05172413793103448
    "P4" = 0.6724137931034483
    "S4" = 0.2758620689655172
    "B6" = 0.05485232067510549
    "D6" = 0.01265822784810127
    "F6" = 0.0379746835443038
    "J6" = 0.3164556962025317
    "O6" = 0.02109704641350211
    "Q6" = 0.1392405063291139
    "R6" = 0.05907172995780591
    "S6" = 0.3586497890295359
    "B7" = 0.1352459016393443
    "D7" = 0.03688524590163934
    "F7" = 0.02868852459016394
    "J7" = 0.180327868852459
    "Q7" = 0.1844262295081967
    "R7" = 0.0778688524590164
    "S7" = 0.3565573770491803
    "B8" = 0.09684684684684684
    "D8" = 0.01351351351351351
    "F8" = 0.06081081081081081
    "J8" = 0.1373873873873874
    "O8" = 0.01126126126126126
    "Q8" = 0.1711711711711712
    "R8" = 0.08108108108108109
    "S8" = 0.4279279279279279
    "B9" = 0.07199999999999999
    "D9" = 0.016
    "F9" = 0.068
    "J9" = 0.12
    "O9" = 0.016
    "Q9" = 0.232
    "R9" = 0.112
    "S9" = 0.364
    "B10" = 0.1404421326397919
    "D10" = 0.02535760728218465
    "F10" = 0.0611183355006502
    "J10" = 0.1287386215864759
    "O10" = 0.00975292587776333
    "Q10" = 0.2106631989596879
    "R10" = 0.06827048114434331
    "S10" = 0.3556566970091027
    "G11" = 0.1173333333333333
    "J11" = 0.1066666666666667
    "K11" = 0.176
    "L11" = 0.5866666666666667
    "S11" = 0.01333333333333333
    "G12" = 0.7488789237668162
    "J12" = 0.2242152466367713
    "K12" = 0.004484304932735426
    "L12" = 0.008968609865470852
    "S12" = 0.01345291479820628
    "G13" = 0.7222222222222222
    "J13" = 0.2407407407407407
    "S13" = 0.03703703703703703
    "F15" = 0.01492537313432836
    "H15" = 0.1940298507462687
    "I15" = 0.07960199004975124
    "J15" = 0.4029850746268657
    "K15" = 0.04477611940298507
    "M15" = 0.004975124378109453
    "O15" = 0.03482587064676617
    "S15" = 0.2238805970149254
    "F16" = 0.02489626556016597
    "H16" = 0.1618257261410788
    "I16" = 0.07468879668049792
    "J16" = 0.4066390041493776
    "K16" = 0.1161825726141079
    "M16" = 0.01244813278008299
    "N16" = 0.004149377593360996
    "O16" = 0.03319502074688797
    "S16" = 0.1659751037344398
    "F17" = 0.01129943502824859
    "H17" = 0.1789077212806026
    "I17" = 0.09416195856873823
    "J17" = 0.4067796610169492
    "K17" = 0.1148775894538606
    "M17" = 0.02824858757062147
    "N17" = 0.001883239171374765
    "O17" = 0.05649717514124294
    "S17" = 0.1073446327683616
    "F18" = 0.025
    "H18" = 0.155
    "I18" = 0.125
    "J18" = 0.44
    "M18" = 0.01
    "O18" = 0.045
    "S18" = 0.1
    "F19" = 0.01970443349753695
    "H19" = 0.1724137931034483
    "I19" = 0.09992962702322308
    "J19" = 0.3856439127375088
    "K19" = 0.1308937368050669
    "M19" = 0.02251935256861365
    "O19" = 0.06403940886699508
    "S19" = 0.1048557353976073
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
